# Generate Report for Handoff
#
# A fresh handoff xliff was generated for
# "68fb0cff-6403-4761-aef1-a62c514756a0.md", so the localization-status
# report's cached timestamps for that row are refreshed on every tab that
# tracks it:
#   - Overview!G6            "Latest HO Xliff Generate Date"
#   - zh-cn!H6                "Latest Handoff Datetime"
#   - de-de!H6                "Latest Handoff Datetime"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G6").Value = "2016-08-12 22:53:05"

$zhCn = $wb.Worksheets.Item("zh-cn")
$zhCn.Range("H6").Value = "2016-08-12 22:52:54"

$deDe = $wb.Worksheets.Item("de-de")
$deDe.Range("H6").Value = "2016-08-12 22:53:05"
